$wb = $excel.ActiveWorkbook

# --- Remove the "scATACseq" entry from the "assay_type list" sheet. ---
$wsAssayType = $wb.Worksheets.Item("assay_type list")
$wsAssayType.Range("A2").EntireRow.Delete()

# --- Remove the "scATACseq" entry from the "transposition_method list" sheet. ---
$wsTransMethod = $wb.Worksheets.Item("transposition_method list")
$wsTransMethod.Range("A2").EntireRow.Delete()

# --- Update the data validation on the main sheet to reflect the shrunk lists. ---
$wsMain = $wb.Worksheets.Item("Export as TSV")

$valAssayType = $wsMain.Range("L2:L1048576").Validation
$valAssayType.Modify(3, 1, 1, "'assay_type list'!`$A`$1:`$A`$3")
$valAssayType.ErrorTitle = "Value must come from list"
$valAssayType.ErrorMessage = "Value must be one of: SNARE-seq2 / sciATACseq / snATACseq."

$valTransMethod = $wsMain.Range("Z2:Z1048576").Validation
$valTransMethod.Modify(3, 1, 1, "'transposition_method list'!`$A`$1:`$A`$4")
$valTransMethod.ErrorTitle = "Value must come from list"
$valTransMethod.ErrorMessage = "Value must be one of: SNARE-Seq2-AC / bulkATACseq / snATACseq / sciATACseq."
